$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LJ Speech")

# Row 2
$ws.Range("B2").Value = "<dere>"

# Row 3
$ws.Range("B3").Value = "<eight>"
$ws.Range("C3").Value = 52

# Row 4
$ws.Range("B4").Value = "<she>"
$ws.Range("C4").Value = 49

# Row 5
$ws.Range("B5").Value = "<out>"
$ws.Range("C5").Value = 55

# Row 6
$ws.Range("B6").Value = "<it>"
$ws.Range("C6").Value = 52

# Row 7
$ws.Range("C7").Value = 51

# Row 8
$ws.Range("C8").Value = 45

# Row 9
$ws.Range("B9").Value = "<then>"
$ws.Range("C9").Value = 15
